$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 122-145: column F (t_match) changes from 7401 to 7400
for ($r = 122; $r -le 145; $r++) {
    $ws.Cells.Item($r, 6).Value = 7400
}

# Rows 242-289: column G (t_b_match) and column H (return) updated for trimmed categories
$ws.Cells.Item(242, 7).Value = 3519.0
$ws.Cells.Item(242, 8).Value = 0.00827862729496309
$ws.Cells.Item(243, 7).Value = 3519.0
$ws.Cells.Item(243, 8).Value = 0.01655725458992618
$ws.Cells.Item(244, 7).Value = 3519.0
$ws.Cells.Item(244, 8).Value = 0.03311450917985236
$ws.Cells.Item(245, 7).Value = 2362.0
$ws.Cells.Item(245, 8).Value = 0.015267799358693662
$ws.Cells.Item(246, 7).Value = 2362.0
$ws.Cells.Item(246, 8).Value = 0.030535598717387324
$ws.Cells.Item(247, 7).Value = 2362.0
$ws.Cells.Item(247, 8).Value = 0.06107119743477465
$ws.Cells.Item(248, 7).Value = 1498.0
$ws.Cells.Item(248, 8).Value = 0.02937290595067568
$ws.Cells.Item(249, 7).Value = 1498.0
$ws.Cells.Item(249, 8).Value = 0.05874581190135136
$ws.Cells.Item(250, 7).Value = 1498.0
$ws.Cells.Item(250, 8).Value = 0.11749162380270271
$ws.Cells.Item(251, 7).Value = 928.0
$ws.Cells.Item(251, 8).Value = 0.056125233303089724
$ws.Cells.Item(252, 7).Value = 928.0
$ws.Cells.Item(252, 8).Value = 0.11225046660617945
$ws.Cells.Item(253, 7).Value = 928.0
$ws.Cells.Item(253, 8).Value = 0.2245009332123589
$ws.Cells.Item(254, 7).Value = 3517.0
$ws.Cells.Item(254, 8).Value = 0.008142532211238385
$ws.Cells.Item(255, 7).Value = 3517.0
$ws.Cells.Item(255, 8).Value = 0.01628506442247677
$ws.Cells.Item(256, 7).Value = 3517.0
$ws.Cells.Item(256, 8).Value = 0.03257012884495354
$ws.Cells.Item(257, 7).Value = 2357.0
$ws.Cells.Item(257, 8).Value = 0.015136305845627345
$ws.Cells.Item(258, 7).Value = 2357.0
$ws.Cells.Item(258, 8).Value = 0.03027261169125469
$ws.Cells.Item(259, 7).Value = 2357.0
$ws.Cells.Item(259, 8).Value = 0.06054522338250938
$ws.Cells.Item(260, 7).Value = 1492.0
$ws.Cells.Item(260, 8).Value = 0.02954830375047737
$ws.Cells.Item(261, 7).Value = 1492.0
$ws.Cells.Item(261, 8).Value = 0.05909660750095474
$ws.Cells.Item(262, 7).Value = 1492.0
$ws.Cells.Item(262, 8).Value = 0.11819321500190948
$ws.Cells.Item(263, 7).Value = 920.0
$ws.Cells.Item(263, 8).Value = 0.05589297240975632
$ws.Cells.Item(264, 7).Value = 920.0
$ws.Cells.Item(264, 8).Value = 0.11178594481951264
$ws.Cells.Item(265, 7).Value = 920.0
$ws.Cells.Item(265, 8).Value = 0.22357188963902527
$ws.Cells.Item(266, 7).Value = 24.0
$ws.Cells.Item(266, 8).Value = -0.012860747149689161
$ws.Cells.Item(267, 7).Value = 24.0
$ws.Cells.Item(267, 8).Value = -0.025721494299378322
$ws.Cells.Item(268, 7).Value = 24.0
$ws.Cells.Item(268, 8).Value = -0.051442988598756643
$ws.Cells.Item(269, 7).Value = 19.0
$ws.Cells.Item(269, 8).Value = -0.014474951449477622
$ws.Cells.Item(270, 7).Value = 19.0
$ws.Cells.Item(270, 8).Value = -0.028949902898955243
$ws.Cells.Item(271, 7).Value = 19.0
$ws.Cells.Item(271, 8).Value = -0.057899805797910486
$ws.Cells.Item(272, 8).Value = -0.009565693984065909
$ws.Cells.Item(273, 8).Value = -0.019131387968131817
$ws.Cells.Item(274, 8).Value = -0.038262775936263635
$ws.Cells.Item(275, 8).Value = 0.011639992939493651
$ws.Cells.Item(276, 8).Value = 0.023279985878987303
$ws.Cells.Item(277, 8).Value = 0.046559971757974605
$ws.Cells.Item(278, 7).Value = 25.0
$ws.Cells.Item(278, 8).Value = -0.013228729039230034
$ws.Cells.Item(279, 7).Value = 25.0
$ws.Cells.Item(279, 8).Value = -0.026457458078460068
$ws.Cells.Item(280, 7).Value = 25.0
$ws.Cells.Item(280, 8).Value = -0.052914916156920136
$ws.Cells.Item(281, 7).Value = 19.0
$ws.Cells.Item(281, 8).Value = -0.01489997528649194
$ws.Cells.Item(282, 7).Value = 19.0
$ws.Cells.Item(282, 8).Value = -0.02979995057298388
$ws.Cells.Item(283, 7).Value = 19.0
$ws.Cells.Item(283, 8).Value = -0.05959990114596776
$ws.Cells.Item(284, 8).Value = -0.009852302152832365
$ws.Cells.Item(285, 8).Value = -0.01970460430566473
$ws.Cells.Item(286, 8).Value = -0.03940920861132946
$ws.Cells.Item(287, 7).Value = 5.0
$ws.Cells.Item(287, 8).Value = 0.02471533791033746
$ws.Cells.Item(288, 7).Value = 5.0
$ws.Cells.Item(288, 8).Value = 0.04943067582067492
$ws.Cells.Item(289, 7).Value = 5.0
$ws.Cells.Item(289, 8).Value = 0.09886135164134983
